# Auto-generated Excel COM-interop script
# Applies numeric cell updates (market-price refresh) across the 8 job sheets
# of the Louisoix_Profits workbook, per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Sheets("ALC")
$ws.Range("H21").Value = 17058.084
$ws.Range("I21").Value = 13959.4
$ws.Range("K21").Value = 13959.4
$ws.Range("M21").Value = -13491.4
$ws.Range("H23").Value = 17058.084
$ws.Range("I23").Value = 13959.4
$ws.Range("K23").Value = 13959.4
$ws.Range("M23").Value = -13725.4
$ws.Range("H106").Value = 4846.8335
$ws.Range("I106").Value = 4846.8335
$ws.Range("K106").Value = 4846.8335
$ws.Range("M106").Value = -4215.8335
$ws.Range("H121").Value = 350
$ws.Range("J121").Value = 350
$ws.Range("L121").Value = 1050
$ws.Range("N121").Value = -4544
$ws.Range("H132").Value = 4759.244
$ws.Range("I132").Value = 2707.162
$ws.Range("K132").Value = 8121.485999999999
$ws.Range("M132").Value = -5591.485999999999
$ws.Range("H137").Value = 7028.1
$ws.Range("I137").Value = 1233
$ws.Range("J137").Value = 9511.714
$ws.Range("K137").Value = 3699
$ws.Range("L137").Value = 28535.142
$ws.Range("M137").Value = -1149
$ws.Range("N137").Value = -33635.142

# --- Sheet: ARM ---
$ws = $wb.Sheets("ARM")
$ws.Range("H4").Value = 588.6
$ws.Range("J4").Value = 748.6667
$ws.Range("L4").Value = 748.6667
$ws.Range("N4").Value = -980.6667
$ws.Range("H61").Value = 1789.7273
$ws.Range("I61").Value = 1436
$ws.Range("J61").Value = 2733
$ws.Range("K61").Value = 1436
$ws.Range("L61").Value = 2733
$ws.Range("M61").Value = -1224
$ws.Range("N61").Value = -3157
$ws.Range("H102").Value = 5218.5864
$ws.Range("I102").Value = 5683.522
$ws.Range("K102").Value = 5683.522
$ws.Range("M102").Value = -4061.522
$ws.Range("H136").Value = 1789.7273
$ws.Range("I136").Value = 1436
$ws.Range("J136").Value = 2733
$ws.Range("K136").Value = 4308
$ws.Range("L136").Value = 8199
$ws.Range("M136").Value = -1758
$ws.Range("N136").Value = -13299

# --- Sheet: BSM ---
$ws = $wb.Sheets("BSM")
$ws.Range("H37").Value = 642.6667
$ws.Range("J37").Value = 414
$ws.Range("L37").Value = 414
$ws.Range("N37").Value = -688
$ws.Range("H102").Value = 11914.625
$ws.Range("I102").Value = 11914.625
$ws.Range("K102").Value = 11914.625
$ws.Range("M102").Value = -8669.625
$ws.Range("H107").Value = 973.55554
$ws.Range("I107").Value = 973.55554
$ws.Range("K107").Value = 973.55554
$ws.Range("M107").Value = 946.44446
$ws.Range("H134").Value = 3634.6667
$ws.Range("I134").Value = 3253.25
$ws.Range("J134").Value = 3939.8
$ws.Range("K134").Value = 9759.75
$ws.Range("L134").Value = 11819.4
$ws.Range("M134").Value = -7224.75
$ws.Range("N134").Value = -16889.4

# --- Sheet: CRP ---
$ws = $wb.Sheets("CRP")
$ws.Range("H7").Value = 154.4
$ws.Range("I7").Value = 166.90909
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 166.90909
$ws.Range("L7").Value = 120
$ws.Range("M7").Value = -53.90908999999999
$ws.Range("N7").Value = -346
$ws.Range("H31").Value = 4413.25
$ws.Range("I31").Value = 2516.3333
$ws.Range("K31").Value = 2516.3333
$ws.Range("M31").Value = -2221.3333
$ws.Range("H34").Value = 4413.25
$ws.Range("I34").Value = 2516.3333
$ws.Range("K34").Value = 2516.3333
$ws.Range("M34").Value = -2314.3333
$ws.Range("H62").Value = 3102
$ws.Range("I62").Value = 2877.5
$ws.Range("K62").Value = 2877.5
$ws.Range("M62").Value = -2253.5
$ws.Range("H65").Value = 3102
$ws.Range("I65").Value = 2877.5
$ws.Range("K65").Value = 14387.5
$ws.Range("M65").Value = -11267.5
$ws.Range("H86").Value = 28882.117
$ws.Range("J86").Value = 5571.9
$ws.Range("L86").Value = 5571.9
$ws.Range("N86").Value = -7817.9
$ws.Range("H89").Value = 28882.117
$ws.Range("J89").Value = 5571.9
$ws.Range("L89").Value = 27859.5
$ws.Range("N89").Value = -39091.5
$ws.Range("H99").Value = 2254.5186
$ws.Range("I99").Value = 1768.125
$ws.Range("J99").Value = 2962
$ws.Range("K99").Value = 1768.125
$ws.Range("L99").Value = 2962
$ws.Range("M99").Value = -270.125
$ws.Range("N99").Value = -5958
$ws.Range("H105").Value = 1626.2941
$ws.Range("I105").Value = 1492.2307
$ws.Range("K105").Value = 1492.2307
$ws.Range("M105").Value = 254.7692999999999
$ws.Range("H126").Value = 2254.5186
$ws.Range("I126").Value = 1768.125
$ws.Range("J126").Value = 2962
$ws.Range("K126").Value = 5304.375
$ws.Range("L126").Value = 8886
$ws.Range("M126").Value = -2834.375
$ws.Range("N126").Value = -13826
$ws.Range("H132").Value = 5006.857
$ws.Range("J132").Value = 10000
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -35060
$ws.Range("H134").Value = 95536.17999999999
$ws.Range("I134").Value = 116099.78
$ws.Range("K134").Value = 348299.34
$ws.Range("M134").Value = -345764.34

# --- Sheet: CUL ---
$ws = $wb.Sheets("CUL")
$ws.Range("H34").Value = 2062.875
$ws.Range("J34").Value = 2100.5
$ws.Range("L34").Value = 6301.5
$ws.Range("N34").Value = -6469.5
$ws.Range("H39").Value = 1577
$ws.Range("J39").Value = 3004
$ws.Range("L39").Value = 9012
$ws.Range("N39").Value = -9600
$ws.Range("H55").Value = 6934
$ws.Range("J55").Value = 8300.799999999999
$ws.Range("L55").Value = 24902.4
$ws.Range("N55").Value = -25256.4
$ws.Range("H92").Value = 486.05
$ws.Range("J92").Value = 562
$ws.Range("L92").Value = 1686
$ws.Range("N92").Value = -4182
$ws.Range("H122").Value = 482.6
$ws.Range("I122").Value = 377.16666
$ws.Range("J122").Value = 640.75
$ws.Range("K122").Value = 3394.49994
$ws.Range("L122").Value = 5766.75
$ws.Range("M122").Value = -944.4999399999997
$ws.Range("N122").Value = -10666.75
$ws.Range("H140").Value = 2524.2
$ws.Range("I140").Value = 2524.2
$ws.Range("K140").Value = 7572.599999999999
$ws.Range("M140").Value = -2392.599999999999

# --- Sheet: GSM ---
$ws = $wb.Sheets("GSM")
$ws.Range("H2").Value = 121.31579
$ws.Range("I2").Value = 139.38461
$ws.Range("J2").Value = 82.166664
$ws.Range("K2").Value = 139.38461
$ws.Range("L2").Value = 82.166664
$ws.Range("M2").Value = -26.38461000000001
$ws.Range("N2").Value = -308.166664
$ws.Range("H105").Value = 74999.2
$ws.Range("J105").Value = 74999.2
$ws.Range("L105").Value = 74999.2
$ws.Range("N105").Value = -81987.2
$ws.Range("H132").Value = 206889.8
$ws.Range("I132").Value = 501224.5
$ws.Range("K132").Value = 1503673.5
$ws.Range("M132").Value = -1501143.5

# --- Sheet: LTW ---
$ws = $wb.Sheets("LTW")
$ws.Range("H7").Value = 8318.214
$ws.Range("I7").Value = 13425
$ws.Range("K7").Value = 13425
$ws.Range("M7").Value = -13313
$ws.Range("H22").Value = 24975.83
$ws.Range("I22").Value = 34246.062
$ws.Range("J22").Value = 3124.5715
$ws.Range("K22").Value = 34246.062
$ws.Range("L22").Value = 3124.5715
$ws.Range("M22").Value = -33951.062
$ws.Range("N22").Value = -3714.5715
$ws.Range("H27").Value = 24975.83
$ws.Range("I27").Value = 34246.062
$ws.Range("J27").Value = 3124.5715
$ws.Range("K27").Value = 34246.062
$ws.Range("L27").Value = 3124.5715
$ws.Range("M27").Value = -34139.062
$ws.Range("N27").Value = -3338.5715
$ws.Range("H40").Value = 3334.2727
$ws.Range("I40").Value = 2641.8
$ws.Range("J40").Value = 3911.3333
$ws.Range("K40").Value = 2641.8
$ws.Range("L40").Value = 3911.3333
$ws.Range("M40").Value = -2505.8
$ws.Range("N40").Value = -4183.3333
$ws.Range("H46").Value = 21359.824
$ws.Range("I46").Value = 32352.3
$ws.Range("K46").Value = 32352.3
$ws.Range("M46").Value = -32164.3
$ws.Range("H101").Value = 23499.5
$ws.Range("J101").Value = 23499.5
$ws.Range("L101").Value = 23499.5
$ws.Range("N101").Value = -29989.5
$ws.Range("H122").Value = 4063.0977
$ws.Range("I122").Value = 3462.0527
$ws.Range("J122").Value = 4582.1816
$ws.Range("K122").Value = 10386.1581
$ws.Range("L122").Value = 13746.5448
$ws.Range("M122").Value = -7936.158100000001
$ws.Range("N122").Value = -18646.5448
$ws.Range("H126").Value = 8318.214
$ws.Range("I126").Value = 13425
$ws.Range("K126").Value = 40275
$ws.Range("M126").Value = -37805
$ws.Range("H132").Value = 63515
$ws.Range("I132").Value = 95530.766
$ws.Range("K132").Value = 286592.298
$ws.Range("M132").Value = -284062.298

# --- Sheet: WVR ---
$ws = $wb.Sheets("WVR")
$ws.Range("H62").Value = 98949.836
$ws.Range("I62").Value = 4345.5713
$ws.Range("J62").Value = 231395.8
$ws.Range("K62").Value = 4345.5713
$ws.Range("L62").Value = 231395.8
$ws.Range("M62").Value = -3721.5713
$ws.Range("N62").Value = -232643.8
$ws.Range("H65").Value = 98949.836
$ws.Range("I65").Value = 4345.5713
$ws.Range("J65").Value = 231395.8
$ws.Range("K65").Value = 21727.8565
$ws.Range("L65").Value = 1156979
$ws.Range("M65").Value = -18607.8565
$ws.Range("N65").Value = -1163219
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H112").Value = 18859.4
$ws.Range("J112").Value = 18859.4
$ws.Range("L112").Value = 18859.4
$ws.Range("N112").Value = -21813.4
$ws.Range("H131").Value = 39996.668
$ws.Range("J131").Value = 39996.668
$ws.Range("L131").Value = 39996.668
$ws.Range("N131").Value = -50076.668
